$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A ("EmpId") so that the existing
# Employee/Hours columns (and their header formatting) shift right:
# Employee -> B, Hours -> C, and the formerly-empty styled cell C12 -> D12.
$ws.Columns.Item(1).Insert()

# New header in A1 ("EmpId"), bold to match the look of the other headers.
$ws.Range("A1").Value = "EmpId"
$ws.Range("A1").Font.Bold = $true

# New "Rate" column header in D1, bold as well.
$ws.Range("D1").Value = "Rate"
$ws.Range("D1").Font.Bold = $true

# Re-assert the existing headers (Employee/Hours kept their formatting
# automatically from the column insert).
$ws.Range("B1").Value = "Employee"
$ws.Range("C1").Value = "Hours"

# EmpId data (row numbers per employee)
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4

# Rate data
$ws.Range("D2").Value = 20
$ws.Range("D3").Value = 20
$ws.Range("D4").Value = 14
$ws.Range("D5").Value = 12

# Column widths for the new/resized columns.
$ws.Columns.Item(1).ColumnWidth = 15.3
$ws.Columns.Item(2).ColumnWidth = 14.59

# Row heights for the header + data rows.
$ws.Rows.Item(1).RowHeight = 13.8
$ws.Rows.Item(2).RowHeight = 13.8
$ws.Rows.Item(3).RowHeight = 13.8
$ws.Rows.Item(4).RowHeight = 13.8
$ws.Rows.Item(5).RowHeight = 13.8

$ws.Range("B9").Select()
